$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 6636.1113
$ws.Range("I40").Value = 7818.5713
$ws.Range("K40").Value = 7818.5713
$ws.Range("M40").Value = -7643.5713
$ws.Range("H69").Value = 7374.75
$ws.Range("I69").Value = 4499
$ws.Range("K69").Value = 13497
$ws.Range("M69").Value = -12623
$ws.Range("H72").Value = 7374.75
$ws.Range("I72").Value = 4499
$ws.Range("K72").Value = 40491
$ws.Range("M72").Value = -36123
$ws.Range("H94").Value = 759.5454999999999
$ws.Range("I94").Value = 759.5454999999999
$ws.Range("K94").Value = 759.5454999999999
$ws.Range("M94").Value = -308.5454999999999
$ws.Range("H112").Value = 37416.934
$ws.Range("J112").Value = 39886.285
$ws.Range("L112").Value = 119658.855
$ws.Range("N112").Value = -121874.855
$ws.Range("H115").Value = 271.875
$ws.Range("I115").Value = 271.875
$ws.Range("K115").Value = 815.625
$ws.Range("M115").Value = 751.375
$ws.Range("H127").Value = 1177.7222
$ws.Range("I127").Value = 881.93335
$ws.Range("J127").Value = 2656.6667
$ws.Range("K127").Value = 2645.80005
$ws.Range("L127").Value = 7970.000100000001
$ws.Range("M127").Value = 2314.19995
$ws.Range("N127").Value = -17890.0001
$ws.Range("H131").Value = 10485.454
$ws.Range("I131").Value = 1992.375
$ws.Range("J131").Value = 33133.668
$ws.Range("K131").Value = 5977.125
$ws.Range("L131").Value = 99401.00399999999
$ws.Range("M131").Value = -937.125
$ws.Range("N131").Value = -109481.004
$ws.Range("H135").Value = 1292.8889
$ws.Range("J135").Value = 2333.3333
$ws.Range("L135").Value = 20999.9997
$ws.Range("N135").Value = -26069.9997
$ws.Range("H137").Value = 1489
$ws.Range("I137").Value = 1258.9375
$ws.Range("K137").Value = 3776.8125
$ws.Range("M137").Value = -1226.8125
$ws.Range("H138").Value = 6852263.5
$ws.Range("I138").Value = 1409.6666
$ws.Range("J138").Value = 7465773
$ws.Range("K138").Value = 4228.9998
$ws.Range("L138").Value = 22397319
$ws.Range("M138").Value = 911.0002000000004
$ws.Range("N138").Value = -22407599
$ws.Range("H140").Value = 92333
$ws.Range("J140").Value = 92333
$ws.Range("L140").Value = 92333
$ws.Range("N140").Value = -102693

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 6347
$ws.Range("I45").Value = 7146.647
$ws.Range("K45").Value = 7146.647
$ws.Range("M45").Value = -6769.647
$ws.Range("H74").Value = 5041.793
$ws.Range("I74").Value = 1098.826
$ws.Range("J74").Value = 20156.5
$ws.Range("K74").Value = 1098.826
$ws.Range("L74").Value = 20156.5
$ws.Range("M74").Value = -224.826
$ws.Range("N74").Value = -21904.5
$ws.Range("H77").Value = 5041.793
$ws.Range("I77").Value = 1098.826
$ws.Range("J77").Value = 20156.5
$ws.Range("K77").Value = 5494.13
$ws.Range("L77").Value = 100782.5
$ws.Range("M77").Value = -1126.13
$ws.Range("N77").Value = -109518.5
$ws.Range("H132").Value = 2957.4
$ws.Range("I132").Value = 2671.9167
$ws.Range("K132").Value = 8015.750100000001
$ws.Range("M132").Value = -5485.750100000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 5470.95
$ws.Range("I86").Value = 6151.75
$ws.Range("J86").Value = 4449.75
$ws.Range("K86").Value = 6151.75
$ws.Range("L86").Value = 4449.75
$ws.Range("M86").Value = -5028.75
$ws.Range("N86").Value = -6695.75
$ws.Range("H89").Value = 5470.95
$ws.Range("I89").Value = 6151.75
$ws.Range("J89").Value = 4449.75
$ws.Range("K89").Value = 30758.75
$ws.Range("L89").Value = 22248.75
$ws.Range("M89").Value = -25142.75
$ws.Range("N89").Value = -33480.75
$ws.Range("H105").Value = 1674.4889
$ws.Range("I105").Value = 1614.8684
$ws.Range("J105").Value = 1998.1428
$ws.Range("K105").Value = 1614.8684
$ws.Range("L105").Value = 1998.1428
$ws.Range("M105").Value = 132.1315999999999
$ws.Range("N105").Value = -5492.1428
$ws.Range("H134").Value = 2333.923
$ws.Range("I134").Value = 2254.2131
$ws.Range("J134").Value = 3549.5
$ws.Range("K134").Value = 6762.6393
$ws.Range("L134").Value = 10648.5
$ws.Range("M134").Value = -4227.6393
$ws.Range("N134").Value = -15718.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 38882.555
$ws.Range("I31").Value = 56975.723
$ws.Range("J31").Value = 2696.2222
$ws.Range("K31").Value = 56975.723
$ws.Range("L31").Value = 2696.2222
$ws.Range("M31").Value = -56680.723
$ws.Range("N31").Value = -3286.2222
$ws.Range("H34").Value = 38882.555
$ws.Range("I34").Value = 56975.723
$ws.Range("J34").Value = 2696.2222
$ws.Range("K34").Value = 56975.723
$ws.Range("L34").Value = 2696.2222
$ws.Range("M34").Value = -56773.723
$ws.Range("N34").Value = -3100.2222
$ws.Range("H105").Value = 1206
$ws.Range("I105").Value = 1227.4445
$ws.Range("J105").Value = 1141.6666
$ws.Range("K105").Value = 1227.4445
$ws.Range("L105").Value = 1141.6666
$ws.Range("M105").Value = 519.5554999999999
$ws.Range("N105").Value = -4635.6666
$ws.Range("H132").Value = 4148.8096
$ws.Range("I132").Value = 4082.0625
$ws.Range("J132").Value = 4362.4
$ws.Range("K132").Value = 12246.1875
$ws.Range("L132").Value = 13087.2
$ws.Range("M132").Value = -9716.1875
$ws.Range("N132").Value = -18147.2

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 4546259
$ws.Range("I113").Value = 8334005
$ws.Range("J113").Value = 964.2
$ws.Range("K113").Value = 25002015
$ws.Range("L113").Value = 2892.6
$ws.Range("M113").Value = -24999845
$ws.Range("N113").Value = -7232.6
$ws.Range("H130").Value = 815
$ws.Range("J130").Value = 1000
$ws.Range("L130").Value = 3000
$ws.Range("N130").Value = -13040
$ws.Range("H131").Value = 34370.195
$ws.Range("I131").Value = 201246.8
$ws.Range("J131").Value = 2278.5386
$ws.Range("K131").Value = 603740.3999999999
$ws.Range("L131").Value = 6835.6158
$ws.Range("M131").Value = -598700.3999999999
$ws.Range("N131").Value = -16915.6158

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H29").Value = 7
$ws.Range("I29").Value = 7
$ws.Range("K29").Value = 7
$ws.Range("M29").Value = 283
$ws.Range("H57").Value = 5000
$ws.Range("I57").Value = 5000
$ws.Range("K57").Value = 5000
$ws.Range("M57").Value = -4180
$ws.Range("H126").Value = 13074.857
$ws.Range("I126").Value = 15275.529
$ws.Range("K126").Value = 45826.587
$ws.Range("M126").Value = -43356.587
$ws.Range("H132").Value = 3679.261
$ws.Range("I132").Value = 3085.1052
$ws.Range("K132").Value = 9255.3156
$ws.Range("M132").Value = -6725.3156

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 2699.652
$ws.Range("I16").Value = 1909.1428
$ws.Range("K16").Value = 1909.1428
$ws.Range("M16").Value = -1739.1428
$ws.Range("H40").Value = 2980.9688
$ws.Range("I40").Value = 2141.0833
$ws.Range("J40").Value = 5500.625
$ws.Range("K40").Value = 2141.0833
$ws.Range("L40").Value = 5500.625
$ws.Range("M40").Value = -2005.0833
$ws.Range("N40").Value = -5772.625
$ws.Range("H64").Value = 100000
$ws.Range("J64").Value = 100000
$ws.Range("L64").Value = 100000
$ws.Range("N64").Value = -100450
$ws.Range("H67").Value = 100000
$ws.Range("J67").Value = 100000
$ws.Range("L67").Value = 100000
$ws.Range("N67").Value = -101560
$ws.Range("H82").Value = 6796.8887
$ws.Range("I82").Value = 10279.363
$ws.Range("K82").Value = 10279.363
$ws.Range("M82").Value = -9918.362999999999
$ws.Range("H85").Value = 6796.8887
$ws.Range("I85").Value = 10279.363
$ws.Range("K85").Value = 10279.363
$ws.Range("M85").Value = -9031.362999999999
$ws.Range("H122").Value = 4566.95
$ws.Range("I122").Value = 3776
$ws.Range("K122").Value = 11328
$ws.Range("M122").Value = -8878

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H94").Value = 13573.5
$ws.Range("J94").Value = 12431.444
$ws.Range("L94").Value = 12431.444
$ws.Range("N94").Value = -14233.444
$ws.Range("H122").Value = 1715.6428
$ws.Range("I122").Value = 1655.3334
$ws.Range("J122").Value = 2077.5
$ws.Range("K122").Value = 4966.0002
$ws.Range("L122").Value = 6232.5
$ws.Range("M122").Value = -2516.0002
$ws.Range("N122").Value = -11132.5
